# Remove the "colony count" index column (column A): the header text and the
# 1..18 row-index values go away, and the colony-count data that was in
# column B (with its "37/37" header) shifts left into column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()

# Update the active selection to match the author's saved cursor position.
$ws.Range("C6").Select()
